$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while forcing text interpretation (so numeric-
# looking strings such as "1.00" or "0.999" are not silently converted to
# numbers) and without leaving a stray number-format style behind.
function Set-TextValue($Cell, $Value) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Value
    $Cell.ClearFormats()
}

# --- Rows 2-27: only the Price (D) and Volume(1h) (E) columns change ---
$deUpdates = @(
    @{ Row = 2; D = '70.996.85'; E = '  -0.21%  ' }
    @{ Row = 3; D = '3.830.14'; E = '  -0.75%  ' }
    @{ Row = 4; D = '0.999'; E = '  -0.07%  ' }
    @{ Row = 5; D = '706.97'; E = '  +1.36%  ' }
    @{ Row = 6; D = '171.63'; E = '  -1.13%  ' }
    @{ Row = 7; D = '3.827.09'; E = '  -0.84%  ' }
    @{ Row = 8; D = '1.00'; E = '  +0.01%  ' }
    @{ Row = 9; D = '0.524'; E = '  -0.55%  ' }
    @{ Row = 10; D = '0.161'; E = '  -1.39%  ' }
    @{ Row = 11; D = '7.43'; E = '  -0.40%  ' }
    @{ Row = 12; D = '0.458'; E = '  -0.75%  ' }
    @{ Row = 13; D = '0.0000254'; E = '  -1.96%  ' }
    @{ Row = 14; D = '36.52'; E = '  -0.56%  ' }
    @{ Row = 15; D = '4.474.44'; E = '  -0.90%  ' }
    @{ Row = 16; D = '3.791.91'; E = '  -1.95%  ' }
    @{ Row = 17; D = '71.010.10'; E = '  -0.29%  ' }
    @{ Row = 18; D = '7.21'; E = '  -0.68%  ' }
    @{ Row = 19; D = '0.115'; E = '  +0.17%  ' }
    @{ Row = 20; D = '17.36'; E = '  -3.04%  ' }
    @{ Row = 21; D = '495.69'; E = '  +1.63%  ' }
    @{ Row = 22; D = '10.65'; E = '  -4.74%  ' }
    @{ Row = 23; D = '0.733'; E = '  +1.59%  ' }
    @{ Row = 24; D = '85.42'; E = '  +0.96%  ' }
    @{ Row = 25; D = '0.0000145'; E = '  -0.75%  ' }
    @{ Row = 26; D = '10.60'; E = '  +0.66%  ' }
    @{ Row = 27; D = '12.09'; E = '  -2.54%  ' }
)

foreach ($u in $deUpdates) {
    Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    Set-TextValue $ws.Cells.Item($u.Row, 5) $u.E
}

# --- Rows 28-51: a new coin (WrappedeETH) was inserted at row 28, pushing
# every following coin down by one row; the former last row (TheGraph)
# drops off the bottom of the list. Overwrite B:E in place (no physical
# row insert) so the sheet keeps its original A1:E51 dimensions. ---
$fullUpdates = @(
    @{ Row = 28; B = 'WrappedeETH'; C = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'; D = '3.981.90'; E = '  -0.77%  ' }
    @{ Row = 29; B = 'Fetch.AI'; C = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'; D = '2.08'; E = '  -3.02%  ' }
    @{ Row = 30; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.00'; E = '  -0.03%  ' }
    @{ Row = 31; B = 'PancakeSwap'; C = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D = '3.09'; E = '  -0.74%  ' }
    @{ Row = 32; B = 'NEARProtocol'; C = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'; D = '7.40'; E = '  -2.91%  ' }
    @{ Row = 33; B = 'ImmutableX'; C = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'; D = '2.23'; E = '  -3.79%  ' }
    @{ Row = 34; B = 'EthereumClassic'; C = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'; D = '29.34'; E = '  -1.69%  ' }
    @{ Row = 35; B = 'Kaspa'; C = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D = '0.175'; E = '  -3.41%  ' }
    @{ Row = 36; B = 'Aptos'; C = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D = '9.17'; E = '  -1.45%  ' }
    @{ Row = 37; B = 'RenzoRestakedETH'; C = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'; D = '3.799.41'; E = '  -0.30%  ' }
    @{ Row = 38; B = 'Binance-PegBSC-USD'; C = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'; D = '1.00'; E = '  -0.32%  ' }
    @{ Row = 39; B = 'Hedera'; C = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; D = '0.102'; E = '  -1.74%  ' }
    @{ Row = 40; B = 'Mantle'; C = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; D = '1.05'; E = '  +3.78%  ' }
    @{ Row = 41; B = 'Stacks'; C = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'; D = '2.33'; E = '  -2.42%  ' }
    @{ Row = 42; B = 'Filecoin'; C = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; D = '5.98'; E = '  -1.37%  ' }
    @{ Row = 43; B = 'dogwifhat'; C = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D = '3.32'; E = '  -3.71%  ' }
    @{ Row = 44; B = 'USDe'; C = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'; D = '1.00'; E = '  -0.03%  ' }
    @{ Row = 45; B = 'FirstDigitalUSD'; C = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D = '1.00'; E = '  -0.02%  ' }
    @{ Row = 46; B = 'FLOKI'; C = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'; D = '0.000312'; E = '  +1.47%  ' }
    @{ Row = 47; B = 'Monero'; C = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; D = '163.42'; E = '  -0.16%  ' }
    @{ Row = 48; B = 'Bittensor'; C = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D = '429.81'; E = '  +3.92%  ' }
    @{ Row = 49; B = 'OKB'; C = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; D = '48.94'; E = '  +0.36%  ' }
    @{ Row = 50; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '8.75'; E = '  +0.37%  ' }
    @{ Row = 51; B = 'ONDO'; C = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'; D = '1.37'; E = '  -1.69%  ' }
)

foreach ($u in $fullUpdates) {
    Set-TextValue $ws.Cells.Item($u.Row, 2) $u.B
    Set-TextValue $ws.Cells.Item($u.Row, 3) $u.C
    Set-TextValue $ws.Cells.Item($u.Row, 4) $u.D
    Set-TextValue $ws.Cells.Item($u.Row, 5) $u.E
}

Write-Output "Update complete"
